$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("H55").Value = 719.1429000000001
$ws.Range("J55").Value = 909.6
$ws.Range("L55").Value = 909.6
$ws.Range("N55").Value = -1337.6

$ws.Range("H101").Value = 477.16666
$ws.Range("I101").Value = 535.6
$ws.Range("K101").Value = 1606.8
$ws.Range("M101").Value = 15.19999999999982

$ws.Range("H113").Value = 4200
$ws.Range("I113").Value = 3933.3333
$ws.Range("K113").Value = 3933.3333
$ws.Range("M113").Value = -679.3332999999998

$ws.Range("H138").Value = 3580.258
$ws.Range("I138").Value = 3152.2964
$ws.Range("J138").Value = 6469
$ws.Range("K138").Value = 9456.889200000001
$ws.Range("L138").Value = 19407
$ws.Range("M138").Value = -4316.889200000001
$ws.Range("N138").Value = -29687

$ws = $wb.Worksheets.Item(2)
$ws.Range("H3").Value = 1924.5
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").Value = ""

$ws = $wb.Worksheets.Item(3)
$ws.Range("H82").Value = 242539
$ws.Range("I82").Value = 6500
$ws.Range("J82").Value = 399898.34
$ws.Range("K82").Value = 6500
$ws.Range("L82").Value = 399898.34
$ws.Range("M82").Value = -6117
$ws.Range("N82").Value = -400664.34

$ws.Range("H85").Value = 242539
$ws.Range("I85").Value = 6500
$ws.Range("J85").Value = 399898.34
$ws.Range("K85").Value = 6500
$ws.Range("L85").Value = 399898.34
$ws.Range("M85").Value = -5174
$ws.Range("N85").Value = -402550.34

$ws.Range("H105").Value = 2851.077
$ws.Range("I105").Value = 2642.2727
$ws.Range("J105").Value = 3999.5
$ws.Range("K105").Value = 2642.2727
$ws.Range("L105").Value = 3999.5
$ws.Range("M105").Value = -895.2727
$ws.Range("N105").Value = -7493.5

$ws = $wb.Worksheets.Item(4)
$ws.Range("H4").Value = 7994
$ws.Range("J4").Value = 7994
$ws.Range("L4").Value = 7994
$ws.Range("N4").Value = -8218

$ws.Range("H7").Value = 113.15
$ws.Range("I7").Value = 46.466667
$ws.Range("J7").Value = 313.2
$ws.Range("K7").Value = 46.466667
$ws.Range("L7").Value = 313.2
$ws.Range("M7").Value = 66.533333
$ws.Range("N7").Value = -539.2

$ws.Range("H22").Value = 330.83334
$ws.Range("I22").Value = 197.5
$ws.Range("J22").Value = 397.5
$ws.Range("K22").Value = 197.5
$ws.Range("L22").Value = 397.5
$ws.Range("M22").Value = 152.5
$ws.Range("N22").Value = -1097.5

$ws.Range("H31").Value = 4431.5454
$ws.Range("I31").Value = 1753.4
$ws.Range("J31").Value = 6663.3335
$ws.Range("K31").Value = 1753.4
$ws.Range("L31").Value = 6663.3335
$ws.Range("M31").Value = -1458.4
$ws.Range("N31").Value = -7253.3335

$ws.Range("H34").Value = 4431.5454
$ws.Range("I34").Value = 1753.4
$ws.Range("J34").Value = 6663.3335
$ws.Range("K34").Value = 1753.4
$ws.Range("L34").Value = 6663.3335
$ws.Range("M34").Value = -1551.4
$ws.Range("N34").Value = -7067.3335

$ws.Range("H134").Value = 2970.3076
$ws.Range("I134").Value = 2891.4375
$ws.Range("J134").Value = 3096.5
$ws.Range("K134").Value = 8674.3125
$ws.Range("L134").Value = 9289.5
$ws.Range("M134").Value = -6139.3125
$ws.Range("N134").Value = -14359.5

$ws = $wb.Worksheets.Item(5)
$ws.Range("H5").Value = 532.0909
$ws.Range("J5").Value = 524.1667
$ws.Range("L5").Value = 1572.5001
$ws.Range("N5").Value = -1796.5001

$ws.Range("H34").Value = 3500
$ws.Range("J34").Value = 3500
$ws.Range("L34").Value = 10500
$ws.Range("N34").Value = -10668

$ws.Range("H36").Value = 1785.7142
$ws.Range("I36").Value = 1785.7142
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 5357.142599999999
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -5188.142599999999
$ws.Range("N36").Value = ""

$ws.Range("H39").Value = 3593.75
$ws.Range("I39").Value = 1437.5
$ws.Range("J39").Value = 5750
$ws.Range("K39").Value = 4312.5
$ws.Range("L39").Value = 17250
$ws.Range("M39").Value = -4018.5
$ws.Range("N39").Value = -17838

$ws.Range("H55").Value = 3876
$ws.Range("J55").Value = 3876
$ws.Range("L55").Value = 11628
$ws.Range("N55").Value = -11982

$ws.Range("H132").Value = 4333.625
$ws.Range("I132").Value = 4333.625
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 39002.625
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -36472.625
$ws.Range("N132").Value = ""

$ws.Range("H135").Value = 532.0909
$ws.Range("J135").Value = 524.1667
$ws.Range("L135").Value = 4717.5003
$ws.Range("N135").Value = -9787.5003

$ws = $wb.Worksheets.Item(6)
$ws.Range("H80").Value = 4944.909
$ws.Range("I80").Value = 4299.6665
$ws.Range("J80").Value = 5186.875
$ws.Range("K80").Value = 4299.6665
$ws.Range("L80").Value = 5186.875
$ws.Range("M80").Value = -3301.6665
$ws.Range("N80").Value = -7182.875

$ws.Range("H83").Value = 4944.909
$ws.Range("I83").Value = 4299.6665
$ws.Range("J83").Value = 5186.875
$ws.Range("K83").Value = 21498.3325
$ws.Range("L83").Value = 25934.375
$ws.Range("M83").Value = -16506.3325
$ws.Range("N83").Value = -35918.375

$ws = $wb.Worksheets.Item(7)
$ws.Range("H22").Value = 6399.875
$ws.Range("I22").Value = 675.25
$ws.Range("J22").Value = 12124.5
$ws.Range("K22").Value = 675.25
$ws.Range("L22").Value = 12124.5
$ws.Range("M22").Value = -380.25
$ws.Range("N22").Value = -12714.5

$ws.Range("H27").Value = 6399.875
$ws.Range("I27").Value = 675.25
$ws.Range("J27").Value = 12124.5
$ws.Range("K27").Value = 675.25
$ws.Range("L27").Value = 12124.5
$ws.Range("M27").Value = -568.25
$ws.Range("N27").Value = -12338.5

$ws.Range("H46").Value = 803.6667
$ws.Range("I46").Value = 650.5
$ws.Range("K46").Value = 650.5
$ws.Range("M46").Value = -462.5

$ws.Range("H55").Value = 426.0625
$ws.Range("I55").Value = 179.90909
$ws.Range("J55").Value = 967.6
$ws.Range("K55").Value = 179.90909
$ws.Range("L55").Value = 967.6
$ws.Range("M55").Value = -6.909089999999992
$ws.Range("N55").Value = -1313.6

$ws = $wb.Worksheets.Item(8)
$ws.Range("H5").Value = 1502496
$ws.Range("J5").Value = 4992
$ws.Range("L5").Value = 4992
$ws.Range("N5").Value = -5216
